$wb = $excel.ActiveWorkbook

# New scrape timestamp (replaces 02:14:53 throughout)
$newTime = "02:37:48"

# --- Sheet 1: LP1912 ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: $newTime"

$ws1.Range("A6").Value = $newTime
$ws1.Range("B6").Value = "02:58"
$ws1.Range("D6").Value = 21

$ws1.Range("A7").Value = $newTime
$ws1.Range("D7").Value = 71

$ws1.Range("A8").Value = $newTime
$ws1.Range("D8").Value = 84

# --- Sheet 2: LP1912-215 ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: $newTime"

$ws2.Range("A6").Value = $newTime
$ws2.Range("B6").Value = "02:58"
$ws2.Range("D6").Value = 21

# --- Sheet 3: 6203-6173 ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: $newTime"
